# Rapport: presentation resultats et discussion
#
# Cell C7 on the single sheet ("Feuil1") held the label for the
# "left-arm perpendicular to trunk" configuration, but with a typo
# ("⏊e" instead of "⏊"). Fix it so it reads the same as C5's label,
# and leave the active selection on C7 (as in the authored edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "Bras gauche `n⏊ au tronc"

$ws.Range("C7").Select()
